$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.125.21"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "2.520.99"
$ws.Range("E3").Value = "  -0.11%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Value = "137.68"
$ws.Range("E6").Value = "  -1.47%  "

$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").Value = "  +0.70%  "

$ws.Range("D9").Value = "2.518.64"
$ws.Range("E9").Value = "  -0.43%  "

$ws.Range("E10").Value = "  +0.12%  "

$ws.Range("E11").Value = "  -2.15%  "

$ws.Range("E12").Value = "  -1.80%  "

$ws.Range("E13").Value = "  -2.40%  "

$ws.Range("D14").Value = "2.967.80"
$ws.Range("E14").Value = "  -0.04%  "

$ws.Range("D15").Value = "23.07"
$ws.Range("E15").Value = "  -1.97%  "

$ws.Range("D16").Value = "58.894.95"
$ws.Range("E16").Value = "  -0.63%  "

$ws.Range("E17").Value = "  -1.61%  "

$ws.Range("D18").Value = "2.519.20"
$ws.Range("E18").Value = "  -0.16%  "

$ws.Range("E19").Value = "  +0.36%  "

$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -0.86%  "

$ws.Range("D21").Value = "325.69"
$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("E22").Value = "  -0.10%  "

$ws.Range("E23").Value = "  +2.09%  "

$ws.Range("D24").Value = "65.81"
$ws.Range("E24").Value = "  +3.46%  "

$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("E27").Value = "  +0.11%  "

$ws.Range("E28").Value = "  -2.46%  "

$ws.Range("D29").Value = "6.69"
$ws.Range("E29").Value = "  -3.33%  "

$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  -1.29%  "

$ws.Range("E31").Value = "  -0.23%  "

$ws.Range("E32").Value = "  +5.65%  "

$ws.Range("D33").Value = "161.99"
$ws.Range("E33").Value = "  -2.02%  "

$ws.Range("E34").Value = "  +0.26%  "

$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("D36").Value = "18.48"
$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("D37").Value = "4.12"
$ws.Range("E37").Value = "  -3.71%  "

$ws.Range("E38").Value = "  -2.31%  "

$ws.Range("D39").Value = "36.58"
$ws.Range("E39").Value = "  -1.08%  "

$ws.Range("D40").Value = "0.819"
$ws.Range("E40").Value = "  +0.09%  "

$ws.Range("E41").Value = "  -1.62%  "

$ws.Range("D42").Value = "286.43"
$ws.Range("E42").Value = "  +2.05%  "

$ws.Range("D43").Value = "5.20"
$ws.Range("E43").Value = "  -1.72%  "

$ws.Range("D44").Value = "132.38"
$ws.Range("E44").Value = "  +7.96%  "

$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  +0.04%  "

$ws.Range("D46").Value = "0.607"
$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("D47").Value = "10.87"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("E48").Value = "  -0.49%  "

$ws.Range("E49").Value = "  -0.97%  "

$ws.Range("E50").Value = "  -1.43%  "

$ws.Range("D51").Value = "17.38"
$ws.Range("E51").Value = "  -2.72%  "
